$d = $word.ActiveDocument

# The original document has a leftover "_GoBack" bookmark (an empty
# range) wrapping the very start of the first paragraph. The edit moves
# this bookmark down to the end of the newly-added last paragraph, so
# first delete it from its current location.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Append the new content after the existing ("git") paragraph: a blank
# paragraph, followed by a paragraph of new text ("Have to get some
# groceries from cnc.---expgit") split across runs the way Word's
# proofing marks ("gramStart"/"gramEnd" around "cnc.") would leave it,
# and finish with a fresh "_GoBack" bookmark at the very end of the
# document (an empty range), mirroring what Word stamps after the last
# edit position.
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p/>
<w:p>
<w:r><w:t xml:space="preserve">Have to get some groceries from </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>cnc.</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>---expgit</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($xml)
